$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell references and their new values, as scraped from the updated
# coinranking.com snapshot (GitHub Actions refresh).
$cells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E11', 'D12', 'E12', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'E20', 'D21', 'E21', 'E22', 'D23', 'E23', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'D30', 'E30', 'D31', 'E31', 'E32', 'D33', 'E33', 'D34', 'E34', 'D35', 'E35', 'D36', 'E36', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'E41', 'D42', 'E42', 'E43', 'E44', 'D45', 'E45', 'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'D51', 'E51')
$newValues = @('28.406.96', '  +1.13%  ', '1.830.87', '  +2.75%  ', '  -0.40%  ', '317.61', '  +0.43%  ', '  -0.36%  ', '0.5332', '  +0.32%  ', '0.4060', '  +8.59%  ', '0.07633', '  +2.66%  ', '41.86', '  +1.69%  ', '6.361', '  +4.62%  ', '  -0.34%  ', '7.603', '  +5.39%  ', '20.89', '  +2.76%  ', '1.828.73', '  +3.85%  ', '0.00001075', '  +2.21%  ', '89.45', '  +0.89%  ', '0.06620', '  +2.16%  ', '  +1.50%  ', '1.001', '  -0.09%  ', '  +3.35%  ', '28.423.30', '  +0.96%  ', '  +1.49%  ', '2.166', '  +3.49%  ', '2.487', '  +9.01%  ', '157.41', '  -0.14%  ', '20.60', '  +1.87%  ', '2.040.71', '  +3.12%  ', '124.50', '  +3.54%  ', '1.126', '  +3.19%  ', '  +5.26%  ', '5.698', '  +3.67%  ', '3.638', '  -0.61%  ', '0.07153', '  +12.60%  ', '0.2262', '  +1.16%  ', '  +3.45%  ', '5.215', '  +5.03%  ', '8.850', '  +4.69%  ', '0.6280', '  +2.13%  ', '  +3.19%  ', '1.189', '  +1.39%  ', '  -0.23%  ', '  -2.29%  ', '13.56', '  +2.24%  ', '3.704', '  +0.96%  ', '0.5863', '  +2.09%  ', '125.82', '  +0.33%  ', '1.991', '  +3.62%  ', '1.205', '  +0.81%  ', '0.06895', '  +0.94%  ')

for ($i = 0; $i -lt $cells.Length; $i++) {
    $cell = $ws.Range($cells[$i])
    # Preserve the cell's existing style, but force a Text number format
    # while assigning the value so Excel doesn't reinterpret numeric-looking
    # strings (e.g. "317.61", "1.001") as actual numbers.
    $origStyle = $cell.Style
    $cell.NumberFormat = '@'
    $cell.Value = $newValues[$i]
    $cell.Style = $origStyle
}
